$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are plain text (e.g. "567.26"); NumberFormat "@" + ClearFormats
# keeps them stored as text (matching the source data) without leaving a numeric
# style behind on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.472.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.433.61"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.428.20"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000177"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.872.73"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.371.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.437.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.19"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.27"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E23").Value = "  +10.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.06"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "616.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.72"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.558.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.49"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  -4.70%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.65"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "145.05"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.81"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.54"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.593"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("E51").Value = "  -1.01%  "
